# CASMNET-858: rename the "INTER_SWITCH_LINKS" tab to "SWITCH_TO_SWITCH" and
# switch which tab/cell is active & selected:
#   - SWITCH_TO_SWITCH (was INTER_SWITCH_LINKS) becomes the selected/active
#     sheet, with cell F41 selected
#   - HARDWARE_MANAGEMENT is no longer the active sheet (its own selection,
#     D28 / A28:XFD28, is left untouched)

$wb = $excel.ActiveWorkbook

# Rename the sheet (keeps its sheetId / position / r:id - just the tab name).
$wsSwitch = $wb.Worksheets.Item("INTER_SWITCH_LINKS")
$wsSwitch.Name = "SWITCH_TO_SWITCH"

# Make it the active tab and move the selection to F41.
$null = $wsSwitch.Activate()
$null = $wsSwitch.Range("F41").Select()
